$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Nombres de prédictions réussies" (B) and "Nombres de matchs" (C)
# for rows 3-11. Column D formulas (% de réussite) recalculate automatically.
$ws.Range("B3").Value = 703
$ws.Range("C3").Value = 1395

$ws.Range("B4").Value = 698
$ws.Range("C4").Value = 1395

$ws.Range("B5").Value = 696
$ws.Range("C5").Value = 1395

$ws.Range("B6").Value = 695
$ws.Range("C6").Value = 1395

$ws.Range("B7").Value = 692
$ws.Range("C7").Value = 1395

$ws.Range("B8").Value = 686
$ws.Range("C8").Value = 1395

$ws.Range("B9").Value = 679
$ws.Range("C9").Value = 1395

$ws.Range("B10").Value = 672
$ws.Range("C10").Value = 1395

$ws.Range("B11").Value = 671
$ws.Range("C11").Value = 1395

# Update the active selection on the sheet (was H4, now G7)
$ws.Range("G7").Select()
